# Update "想去人数" (interest count) values in column F on the "展览"
# and "全部类型" sheets, reflecting the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 318
$ws1.Range("F5").Value = 1287
$ws1.Range("F6").Value = 77
$ws1.Range("F7").Value = 2142
$ws1.Range("F11").Value = 4775
$ws1.Range("F14").Value = 297
$ws1.Range("F21").Value = 3696
$ws1.Range("F22").Value = 496
$ws1.Range("F23").Value = 604
$ws1.Range("F27").Value = 111
$ws1.Range("F34").Value = 831
$ws1.Range("F35").Value = 2307

# Sheet "全部类型" (sheet4): row -> new value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 318
$ws4.Range("F5").Value = 1287
$ws4.Range("F6").Value = 77
$ws4.Range("F7").Value = 2142
$ws4.Range("F11").Value = 4775
$ws4.Range("F14").Value = 297
$ws4.Range("F21").Value = 3696
$ws4.Range("F22").Value = 496
$ws4.Range("F23").Value = 604
$ws4.Range("F27").Value = 111
$ws4.Range("F35").Value = 831
$ws4.Range("F36").Value = 2307
